$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update worker record in row 16 ---
# N° Doc Trabajador
$ws.Range("C16").Value = "1044907764"
# Nombre Trabajador
$ws.Range("D16").Value = "CATRY LUZ AMADOR MUÑIZ"
# Periodo Mora (now centered)
$ws.Range("E16").Value = "2509"
$ws.Range("E16").HorizontalAlignment = -4108
# Valor Mora (table)
$ws.Range("F16").Value = 47450

# --- Update VALOR MORA total at top of sheet ---
$ws.Range("E11").Value = 47450

# --- Column D best-fit width recalculated by Excel for the new name ---
$ws.Columns.Item(4).ColumnWidth = 25.54296875
